$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 446.1111
$ws.Range("I99").Value = 314.375
$ws.Range("K99").Value = 943.125
$ws.Range("M99").Value = 554.875

$ws.Range("H100").Value = 2674.125
$ws.Range("I100").Value = 2350.5
$ws.Range("K100").Value = 2350.5
$ws.Range("M100").Value = -1809.5

$ws.Range("H112").Value = 3515.4285
$ws.Range("J112").Value = 3501.3333
$ws.Range("L112").Value = 10503.9999
$ws.Range("N112").Value = -12719.9999

$ws.Range("H135").Value = 698.25
$ws.Range("I135").Value = 698.25
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 6284.25
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -3749.25

$ws.Range("H138").Value = 2547.2927
$ws.Range("I138").Value = 2893.238
$ws.Range("J138").Value = 2184.05
$ws.Range("K138").Value = 8679.714
$ws.Range("L138").Value = 6552.150000000001
$ws.Range("M138").Value = -3539.714
$ws.Range("N138").Value = -16832.15

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 14801
$ws.Range("J104").Value = 14801
$ws.Range("L104").Value = 14801
$ws.Range("N104").Value = -21789

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2143.1428
$ws.Range("I31").Value = 1640.25
$ws.Range("K31").Value = 1640.25
$ws.Range("M31").Value = -1345.25

$ws.Range("H34").Value = 2143.1428
$ws.Range("I34").Value = 1640.25
$ws.Range("K34").Value = 1640.25
$ws.Range("M34").Value = -1438.25

$ws.Range("H35").Value = 1230
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H107").Value = 759.1429000000001
$ws.Range("I107").Value = 142.75
$ws.Range("J107").Value = 1581
$ws.Range("K107").Value = 142.75
$ws.Range("L107").Value = 1581
$ws.Range("M107").Value = 1777.25
$ws.Range("N107").Value = -5421

$ws.Range("H132").Value = 2218
$ws.Range("I132").Value = 2132
$ws.Range("K132").Value = 6396
$ws.Range("M132").Value = -3866

$ws.Range("H134").Value = 2055.5
$ws.Range("I134").Value = 1740.6666
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 5221.9998
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -2686.9998
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2860.5715
$ws.Range("I5").Value = 2136
$ws.Range("J5").Value = 4672
$ws.Range("K5").Value = 6408
$ws.Range("L5").Value = 14016
$ws.Range("M5").Value = -6296
$ws.Range("N5").Value = -14240

$ws.Range("H39").Value = 2893
$ws.Range("J39").Value = 280
$ws.Range("L39").Value = 840
$ws.Range("N39").Value = -1428

$ws.Range("H41").Value = 1500
$ws.Range("I41").Value = 1500
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 4500
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -4162

$ws.Range("H63").Value = 11292
$ws.Range("I63").Value = 13115
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 39345
$ws.Range("L63").Value = 12000
$ws.Range("M63").Value = -38596
$ws.Range("N63").Value = -13498

$ws.Range("H64").Value = 16507
$ws.Range("I64").Value = 15000
$ws.Range("J64").Value = 18014
$ws.Range("K64").Value = 45000
$ws.Range("L64").Value = 54042
$ws.Range("M64").Value = -44730
$ws.Range("N64").Value = -54582

$ws.Range("H66").Value = 11292
$ws.Range("I66").Value = 13115
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 118035
$ws.Range("L66").Value = 36000
$ws.Range("M66").Value = -114291
$ws.Range("N66").Value = -43488

$ws.Range("H67").Value = 16507
$ws.Range("I67").Value = 15000
$ws.Range("J67").Value = 18014
$ws.Range("K67").Value = 45000
$ws.Range("L67").Value = 54042
$ws.Range("M67").Value = -44064
$ws.Range("N67").Value = -55914

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("N80").Value = 0

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("N83").Value = 0

$ws.Range("H129").Value = 771525.0600000001
$ws.Range("I129").Value = 1108
$ws.Range("J129").Value = 1253035.8
$ws.Range("K129").Value = 3324
$ws.Range("L129").Value = 3759107.4
$ws.Range("M129").Value = 1676
$ws.Range("N129").Value = -3769107.4

$ws.Range("H135").Value = 2860.5715
$ws.Range("I135").Value = 2136
$ws.Range("J135").Value = 4672
$ws.Range("K135").Value = 19224
$ws.Range("L135").Value = 42048
$ws.Range("M135").Value = -16689
$ws.Range("N135").Value = -47118

$ws.Range("H137").Value = 1584
$ws.Range("I137").Value = 1584
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 4752
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = 348

$ws.Range("H139").Value = 1867.2
$ws.Range("I139").Value = 1778.5
$ws.Range("K139").Value = 5335.5
$ws.Range("M139").Value = -195.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2022.7222
$ws.Range("I80").Value = 1616.6666
$ws.Range("J80").Value = 2103.9333
$ws.Range("K80").Value = 1616.6666
$ws.Range("L80").Value = 2103.9333
$ws.Range("M80").Value = -618.6666
$ws.Range("N80").Value = -4099.933300000001

$ws.Range("H83").Value = 2022.7222
$ws.Range("I83").Value = 1616.6666
$ws.Range("J83").Value = 2103.9333
$ws.Range("K83").Value = 8083.333000000001
$ws.Range("L83").Value = 10519.6665
$ws.Range("M83").Value = -3091.333000000001
$ws.Range("N83").Value = -20503.6665

$ws.Range("H102").Value = 1707.5454
$ws.Range("I102").Value = 1628.3
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 1628.3
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -6.299999999999955
$ws.Range("N102").Value = -5744

$ws.Range("H107").Value = 2927.25
$ws.Range("I107").Value = 119.75
$ws.Range("J107").Value = 5734.75
$ws.Range("K107").Value = 119.75
$ws.Range("L107").Value = 5734.75
$ws.Range("M107").Value = 1800.25
$ws.Range("N107").Value = -9574.75

$ws.Range("H122").Value = 1638.8572
$ws.Range("I122").Value = 1702.6666
$ws.Range("J122").Value = 1524
$ws.Range("K122").Value = 5107.9998
$ws.Range("L122").Value = 4572
$ws.Range("M122").Value = -2657.9998
$ws.Range("N122").Value = -9472

$ws.Range("H126").Value = 3512.7144
$ws.Range("I126").Value = 3549
$ws.Range("J126").Value = 3295
$ws.Range("K126").Value = 10647
$ws.Range("L126").Value = 9885
$ws.Range("M126").Value = -8177
$ws.Range("N126").Value = -14825

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 290
$ws.Range("I9").Value = 236.66667
$ws.Range("J9").Value = 450
$ws.Range("K9").Value = 236.66667
$ws.Range("L9").Value = 450
$ws.Range("M9").Value = -12.66667000000001
$ws.Range("N9").Value = -898

$ws.Range("H22").Value = 2792
$ws.Range("I22").Value = 2490
$ws.Range("J22").Value = 4000
$ws.Range("K22").Value = 2490
$ws.Range("L22").Value = 4000
$ws.Range("M22").Value = -2195
$ws.Range("N22").Value = -4590

$ws.Range("H27").Value = 2792
$ws.Range("I27").Value = 2490
$ws.Range("J27").Value = 4000
$ws.Range("K27").Value = 2490
$ws.Range("L27").Value = 4000
$ws.Range("M27").Value = -2383
$ws.Range("N27").Value = -4214

$ws.Range("H40").Value = 2703.1667
$ws.Range("I40").Value = 2174.9375
$ws.Range("K40").Value = 2174.9375
$ws.Range("M40").Value = -2038.9375

$ws.Range("H46").Value = 3250.1667
$ws.Range("I46").Value = 2799.5
$ws.Range("J46").Value = 3475.5
$ws.Range("K46").Value = 2799.5
$ws.Range("L46").Value = 3475.5
$ws.Range("M46").Value = -2611.5
$ws.Range("N46").Value = -3851.5

$ws.Range("H100").Value = 3500
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 3500
$ws.Range("K100").Value = 0
$ws.Range("L100").ClearContents()
$ws.Range("M100").Value = 3500
$ws.Range("N100").Value = -4582

$ws.Range("H111").Value = 40000
$ws.Range("J111").Value = 40000
$ws.Range("L111").Value = 40000
$ws.Range("N111").Value = -48180

$ws.Range("H132").Value = 2555.9285
$ws.Range("I132").Value = 2367.5386
$ws.Range("K132").Value = 7102.6158
$ws.Range("M132").Value = -4572.6158

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 13904.833
$ws.Range("J62").Value = 13161.777
$ws.Range("L62").Value = 13161.777
$ws.Range("N62").Value = -14409.777

$ws.Range("H65").Value = 13904.833
$ws.Range("J65").Value = 13161.777
$ws.Range("L65").Value = 65808.88499999999
$ws.Range("N65").Value = -72048.88499999999

$ws.Range("H132").Value = 3661.389
$ws.Range("I132").Value = 3960.8
$ws.Range("K132").Value = 11882.4
$ws.Range("M132").Value = -9352.400000000001

$ws.Range("H136").Value = 2796.2903
$ws.Range("I136").Value = 2719.72
$ws.Range("J136").Value = 3115.3333
$ws.Range("K136").Value = 8159.16
$ws.Range("L136").Value = 9345.999899999999
$ws.Range("M136").Value = -5609.16
$ws.Range("N136").Value = -14445.9999
